$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.520.09'
$ws.Range("E2").Value = '  +3.68%  '

$ws.Range("D3").Value = '1.587.10'
$ws.Range("E3").Value = '  +0.89%  '

$ws.Range("D5").Value = "'212.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.62%  '

$ws.Range("D6").Value = "'0.493"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("E7").Value = '  +1.00%  '

$ws.Range("D8").Value = "'24.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.21%  '

$ws.Range("E9").Value = '  +0.57%  '

$ws.Range("D10").Value = "'0.0601"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.91%  '

$ws.Range("D11").Value = "'0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.75%  '

$ws.Range("D12").Value = '1.813.91'
$ws.Range("E12").Value = '  +0.91%  '

$ws.Range("D13").Value = '1.602.26'
$ws.Range("E13").Value = '  +2.01%  '

$ws.Range("E14").Value = '  +1.77%  '

$ws.Range("E15").Value = '  -0.53%  '

$ws.Range("D16").Value = '28.533.29'
$ws.Range("E16").Value = '  +3.79%  '

$ws.Range("E17").Value = '  +1.03%  '

$ws.Range("D18").Value = "'231.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.24%  '

$ws.Range("D19").Value = "'7.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("E20").Value = '  +0.18%  '

$ws.Range("E21").Value = '  +0.97%  '

$ws.Range("E22").Value = '  -1.90%  '

$ws.Range("D23").Value = "'9.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.98%  '

$ws.Range("E24").Value = '  +1.97%  '

$ws.Range("D25").Value = "'151.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.99%  '

$ws.Range("D26").Value = "'15.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.34%  '

$ws.Range("E27").Value = '  -0.90%  '

$ws.Range("E28").Value = '  -0.82%  '

$ws.Range("E29").Value = '  +0.97%  '

$ws.Range("E30").Value = '  -1.30%  '

$ws.Range("E31").Value = '  -0.69%  '

$ws.Range("D32").Value = "'3.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.20%  '

$ws.Range("D33").Value = "'3.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.24%  '

$ws.Range("D34").Value = '1.389.83'
$ws.Range("E34").Value = '  -4.70%  '

$ws.Range("E35").Value = '  -1.42%  '

$ws.Range("E36").Value = '  -10.75%  '

$ws.Range("E37").Value = '  +1.13%  '

$ws.Range("E38").Value = '  +10.70%  '

$ws.Range("E39").Value = '  -0.57%  '

$ws.Range("E40").Value = '  +0.16%  '

$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("E42").Value = '  +0.99%  '

$ws.Range("E43").Value = '  -0.53%  '

$ws.Range("E44").Value = '  +0.16%  '

$ws.Range("E45").Value = '  +0.73%  '

$ws.Range("D46").Value = "'62.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.96%  '

$ws.Range("D47").Value = '1.724.71'
$ws.Range("E47").Value = '  +0.91%  '

$ws.Range("D48").Value = "'2.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.42%  '

$ws.Range("D49").Value = "'87.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.04%  '

$ws.Range("E50").Value = '  +0.49%  '

$ws.Range("E51").Value = '  -1.31%  '
